$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: target cell and the new value to write.
# ForceText entries get a leading apostrophe so Excel keeps the
# numeric-looking price strings (e.g. trailing zeros) as literal text
# instead of silently converting them to numbers.
$changes = @(
    @{Cell='D2'; Value='26.791.42'},
    @{Cell='E2'; Value='  -4.11%  '},
    @{Cell='D3'; Value='1.721.78'},
    @{Cell='E3'; Value='  -2.44%  '},
    @{Cell='E4'; Value='  +0.19%  '},
    @{Cell='D5'; Value='''309.53'},
    @{Cell='E5'; Value='  -5.79%  '},
    @{Cell='D6'; Value='''1.003'},
    @{Cell='E6'; Value='  +0.20%  '},
    @{Cell='D7'; Value='''0.4848'},
    @{Cell='E7'; Value='  +3.51%  '},
    @{Cell='D8'; Value='''0.3474'},
    @{Cell='E8'; Value='  -1.38%  '},
    @{Cell='D9'; Value='''42.65'},
    @{Cell='E9'; Value='  -2.35%  '},
    @{Cell='D10'; Value='''0.07223'},
    @{Cell='E10'; Value='  -2.03%  '},
    @{Cell='D11'; Value='''1.048'},
    @{Cell='E11'; Value='  -3.04%  '},
    @{Cell='D12'; Value='''1.003'},
    @{Cell='E12'; Value='  +0.23%  '},
    @{Cell='D13'; Value='''19.76'},
    @{Cell='E13'; Value='  -4.03%  '},
    @{Cell='D14'; Value='''5.856'},
    @{Cell='E14'; Value='  -2.29%  '},
    @{Cell='D15'; Value='1.721.30'},
    @{Cell='E15'; Value='  -2.39%  '},
    @{Cell='D16'; Value='''6.803'},
    @{Cell='E16'; Value='  -5.13%  '},
    @{Cell='D17'; Value='''86.35'},
    @{Cell='E17'; Value='  -6.35%  '},
    @{Cell='E18'; Value='  -1.54%  '},
    @{Cell='D19'; Value='''0.06405'},
    @{Cell='E19'; Value='  -0.24%  '},
    @{Cell='D21'; Value='''16.52'},
    @{Cell='E21'; Value='  -2.23%  '},
    @{Cell='D22'; Value='''5.708'},
    @{Cell='E22'; Value='  -1.09%  '},
    @{Cell='D23'; Value='26.867.84'},
    @{Cell='E23'; Value='  -3.92%  '},
    @{Cell='D24'; Value='''10.90'},
    @{Cell='E24'; Value='  -1.96%  '},
    @{Cell='D25'; Value='''2.056'},
    @{Cell='E25'; Value='  -4.45%  '},
    @{Cell='D26'; Value='''154.70'},
    @{Cell='E26'; Value='  -4.80%  '},
    @{Cell='D27'; Value='''19.82'},
    @{Cell='E27'; Value='  -0.81%  '},
    @{Cell='D28'; Value='1.919.29'},
    @{Cell='E28'; Value='  -2.36%  '},
    @{Cell='E29'; Value='  -5.11%  '},
    @{Cell='D30'; Value='''120.24'},
    @{Cell='E30'; Value='  -2.11%  '},
    @{Cell='E31'; Value='  -3.45%  '},
    @{Cell='D32'; Value='''0.09292'},
    @{Cell='E32'; Value='  -0.01%  '},
    @{Cell='E33'; Value='  -2.16%  '},
    @{Cell='D34'; Value='''5.362'},
    @{Cell='E34'; Value='  -3.22%  '},
    @{Cell='D35'; Value='''0.05905'},
    @{Cell='E35'; Value='  -2.86%  '},
    @{Cell='D36'; Value='''0.02176'},
    @{Cell='E36'; Value='  -3.78%  '},
    @{Cell='D37'; Value='''1.426'},
    @{Cell='E37'; Value='  +0.40%  '},
    @{Cell='D38'; Value='''10.93'},
    @{Cell='E38'; Value='  -6.33%  '},
    @{Cell='E39'; Value='  +0.19%  '},
    @{Cell='D40'; Value='''0.1982'},
    @{Cell='E40'; Value='  -3.92%  '},
    @{Cell='D41'; Value='''4.721'},
    @{Cell='E41'; Value='  -3.68%  '},
    @{Cell='D42'; Value='''0.5957'},
    @{Cell='E42'; Value='  -2.84%  '},
    @{Cell='D43'; Value='''1.117'},
    @{Cell='E43'; Value='  -5.89%  '},
    @{Cell='D44'; Value='''7.430'},
    @{Cell='E44'; Value='  -4.03%  '},
    @{Cell='D45'; Value='''12.82'},
    @{Cell='E45'; Value='  -2.63%  '},
    @{Cell='D46'; Value='''3.573'},
    @{Cell='E46'; Value='  -4.41%  '},
    @{Cell='D47'; Value='''0.5581'},
    @{Cell='E47'; Value='  -3.47%  '},
    @{Cell='D48'; Value='''118.95'},
    @{Cell='E48'; Value='  -3.54%  '},
    @{Cell='D49'; Value='''1.830'},
    @{Cell='E49'; Value='  -5.10%  '},
    @{Cell='B50'; Value='Cronos'},
    @{Cell='C50'; Value='https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'},
    @{Cell='D50'; Value='''0.06620'},
    @{Cell='E50'; Value='  -2.89%  '},
    @{Cell='B51'; Value='EOS'},
    @{Cell='C51'; Value='https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'},
    @{Cell='D51'; Value='''1.091'},
    @{Cell='E51'; Value='  -2.83%  '}
)

foreach ($ch in $changes) {
    $ws.Range($ch.Cell).Value = $ch.Value
}

Write-Host "Applied $($changes.Count) cell updates"